# Resolve the workbook / primary "data" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query timestamps on the "data" sheet (F2:F5) ---
$ws.Range("F2").Value = "2021-10-05 14:35:40.750210"
$ws.Range("F3").Value = "2021-10-05 14:35:40.750218"
$ws.Range("F4").Value = "2021-10-05 14:35:40.750221"
$ws.Range("F5").Value = "2021-10-05 14:35:40.750224"

# --- Add a new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Apply the same (bold / bordered / centered) header style used on "data"!B1
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Sarcoma"
$meta.Range("C2").Value = 3281
# "0.2" must stay a text string, not be coerced to the number 0.2
$meta.Range("D2").Value = "'0.2"
$meta.Range("E2").Value = "2020-08-10T07:04:06.291247Z"
$meta.Range("F2").Value = "2021-10-05 14:35:40.746387"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3281/?format=json"

# Apply the same style as "data"!A2 to metadata!A2 (numeric index cell style)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# D2 should carry no special formatting (plain/default style), same as the
# rest of the data row cells such as "data"!B2 - copy that plain format over
# to drop the transient quote-prefix formatting picked up above.
$ws.Range("B2").Copy()
$meta.Range("D2").PasteSpecial(-4122)

# Keep "data" as the active/selected tab (book view is unchanged by the edit)
$ws.Activate()
